# Horarios actualizados Linea 141 - 1307
# Applies the scrape-refresh diff (update timestamp 20:00:07 -> 20:28:23,
# reshuffled rows, and newly scraped arrival rows) across the three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 20:28:23"
$ws1.Range("A3").Value = "Total filas: 315"

# Rows 15/16 swapped their "Linea" label
$ws1.Cells.Item(15, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(16, 3).Value = "225_GOMEZ"

# Rows 23/24 swapped
$ws1.Cells.Item(23, 1).Value = "06:15:23"
$ws1.Cells.Item(23, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(23, 4).Value = 66

$ws1.Cells.Item(24, 1).Value = "06:46:40"
$ws1.Cells.Item(24, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(24, 4).Value = 35

# Rows 71/72 swapped
$ws1.Cells.Item(71, 1).Value = "08:49:06"
$ws1.Cells.Item(71, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(71, 4).Value = 42

$ws1.Cells.Item(72, 1).Value = "08:14:55"
$ws1.Cells.Item(72, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(72, 4).Value = 77

# Rows 103/104 swapped
$ws1.Cells.Item(103, 1).Value = "10:32:07"
$ws1.Cells.Item(103, 3).Value = "14_ABASTO"
$ws1.Cells.Item(103, 4).Value = 43

$ws1.Cells.Item(104, 1).Value = "11:01:19"
$ws1.Cells.Item(104, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(104, 4).Value = 14

# Rows 117/118 swapped their "Linea" label
$ws1.Cells.Item(117, 3).Value = "17_ROMERO"
$ws1.Cells.Item(118, 3).Value = "16_SANTA ANA"

# Rows 133/134 swapped their "Linea" label
$ws1.Cells.Item(133, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(134, 3).Value = "23_HERNANDEZ"

# Rows 140/141 swapped
$ws1.Cells.Item(140, 1).Value = "12:43:13"
$ws1.Cells.Item(140, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(140, 4).Value = 20

$ws1.Cells.Item(141, 1).Value = "12:18:38"
$ws1.Cells.Item(141, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(141, 4).Value = 45

# Rows 187/188 swapped
$ws1.Cells.Item(187, 1).Value = "14:58:38"
$ws1.Cells.Item(187, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(187, 4).Value = 43

$ws1.Cells.Item(188, 1).Value = "15:34:15"
$ws1.Cells.Item(188, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(188, 4).Value = 7

# Rows 190/192 swapped their "Linea" label
$ws1.Cells.Item(190, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(192, 3).Value = "16_P MOR-SANTA ANA"

# Rows 199/200 swapped
$ws1.Cells.Item(199, 1).Value = "16:02:30"
$ws1.Cells.Item(199, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(199, 4).Value = 3

$ws1.Cells.Item(200, 1).Value = "14:58:38"
$ws1.Cells.Item(200, 3).Value = "14_ABASTO"
$ws1.Cells.Item(200, 4).Value = 67

# Rows 280/281 swapped
$ws1.Cells.Item(280, 1).Value = "19:14:15"
$ws1.Cells.Item(280, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(280, 4).Value = 9

$ws1.Cells.Item(281, 1).Value = "18:01:05"
$ws1.Cells.Item(281, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(281, 4).Value = 82

# Tail rows refreshed with the new scrape timestamp / recomputed ETAs
$ws1.Cells.Item(299, 1).Value = "20:28:23"
$ws1.Cells.Item(299, 4).Value = 3

$ws1.Cells.Item(300, 1).Value = "20:28:23"
$ws1.Cells.Item(300, 4).Value = 6

$ws1.Cells.Item(301, 1).Value = "20:28:23"
$ws1.Cells.Item(301, 4).Value = 18

$ws1.Cells.Item(302, 1).Value = "20:28:23"
$ws1.Cells.Item(302, 4).Value = 19

$ws1.Cells.Item(306, 1).Value = "20:28:23"
$ws1.Cells.Item(306, 4).Value = 27

$ws1.Cells.Item(308, 1).Value = "20:28:23"
$ws1.Cells.Item(308, 4).Value = 29

$ws1.Cells.Item(309, 1).Value = "20:28:23"
$ws1.Cells.Item(309, 4).Value = 39

$ws1.Cells.Item(310, 1).Value = "20:28:23"
$ws1.Cells.Item(310, 4).Value = 42

$ws1.Cells.Item(311, 1).Value = "20:28:23"
$ws1.Cells.Item(311, 4).Value = 60

$ws1.Cells.Item(314, 1).Value = "20:28:23"
$ws1.Cells.Item(314, 4).Value = 66

# Old row 315 (14X44_ABASTO) becomes 16_SANTA ANA with an earlier ETA
$ws1.Cells.Item(315, 1).Value = "20:28:23"
$ws1.Cells.Item(315, 2).Value = "21:34"
$ws1.Cells.Item(315, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(315, 4).Value = 66

# New row 316 (23_HERNANDEZ)
$ws1.Cells.Item(316, 1).Value = "20:28:23"
$ws1.Cells.Item(316, 2).Value = "21:45"
$ws1.Cells.Item(316, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(316, 4).Value = 77
$ws1.Cells.Item(316, 5).Value = "LP1912"

# New row 317 (14X44_ABASTO, the old row-315 content re-scraped)
$ws1.Cells.Item(317, 1).Value = "20:28:23"
$ws1.Cells.Item(317, 2).Value = "21:46"
$ws1.Cells.Item(317, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(317, 4).Value = 78
$ws1.Cells.Item(317, 5).Value = "LP1912"

# Row 318: what used to be row 316 (23_HERNANDEZ / 21:48 / 108), shifted down
$ws1.Cells.Item(318, 1).Value = "20:00:07"
$ws1.Cells.Item(318, 2).Value = "21:48"
$ws1.Cells.Item(318, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(318, 4).Value = 108
$ws1.Cells.Item(318, 5).Value = "LP1912"

# New row 319 (15_ABASTO)
$ws1.Cells.Item(319, 1).Value = "20:28:23"
$ws1.Cells.Item(319, 2).Value = "22:04"
$ws1.Cells.Item(319, 3).Value = "15_ABASTO"
$ws1.Cells.Item(319, 4).Value = 96
$ws1.Cells.Item(319, 5).Value = "LP1912"

# New row 320 (14_ABASTO)
$ws1.Cells.Item(320, 1).Value = "20:28:23"
$ws1.Cells.Item(320, 2).Value = "22:11"
$ws1.Cells.Item(320, 3).Value = "14_ABASTO"
$ws1.Cells.Item(320, 4).Value = 103
$ws1.Cells.Item(320, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 20:28:23"

$ws2.Cells.Item(52, 1).Value = "20:28:23"
$ws2.Cells.Item(52, 4).Value = 19

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 20:28:23"
$ws3.Range("A3").Value = "Total filas: 44"

$ws3.Cells.Item(47, 1).Value = "20:28:23"
$ws3.Cells.Item(47, 4).Value = 24

# New row 48 (215C_LA PLATA)
$ws3.Cells.Item(48, 1).Value = "20:28:23"
$ws3.Cells.Item(48, 2).Value = "21:30"
$ws3.Cells.Item(48, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(48, 4).Value = 62
$ws3.Cells.Item(48, 5).Value = "L6203"

# New row 49 (215B_LP-P MOR-40 Y 115)
$ws3.Cells.Item(49, 1).Value = "20:28:23"
$ws3.Cells.Item(49, 2).Value = "22:20"
$ws3.Cells.Item(49, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(49, 4).Value = 112
$ws3.Cells.Item(49, 5).Value = "L6173"
